# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newer counts, as captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 199
$ws1.Range("F3").Value  = 5465
$ws1.Range("F7").Value  = 637
$ws1.Range("F8").Value  = 607
$ws1.Range("F9").Value  = 1065
$ws1.Range("F11").Value = 1515
$ws1.Range("F12").Value = 4813
$ws1.Range("F14").Value = 212
$ws1.Range("F15").Value = 186
$ws1.Range("F16").Value = 101
$ws1.Range("F17").Value = 3919
$ws1.Range("F19").Value = 1129
$ws1.Range("F20").Value = 112
$ws1.Range("F22").Value = 207
$ws1.Range("F23").Value = 38
$ws1.Range("F24").Value = 146
$ws1.Range("F25").Value = 53
$ws1.Range("F28").Value = 330

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 199
$ws4.Range("F4").Value  = 5465
$ws4.Range("F8").Value  = 637
$ws4.Range("F9").Value  = 607
$ws4.Range("F10").Value = 1065
$ws4.Range("F12").Value = 1515
$ws4.Range("F13").Value = 4813
$ws4.Range("F15").Value = 212
$ws4.Range("F16").Value = 186
$ws4.Range("F17").Value = 101
$ws4.Range("F18").Value = 3919
$ws4.Range("F20").Value = 1129
$ws4.Range("F21").Value = 112
$ws4.Range("F23").Value = 207
$ws4.Range("F24").Value = 38
$ws4.Range("F25").Value = 146
$ws4.Range("F26").Value = 53
$ws4.Range("F29").Value = 330
